# "missed in last commit": add the missing {{renewal}} placeholder cell
# in row 2 (the template/sample-data row) between Certificate Type (O)
# and GroupCvrgLevel (Q). The cell previously existed but was left blank
# (it only carried a stray cell style); give it the same look as its
# neighbours and fill in the placeholder text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make P2 match the formatting already used across the rest of row 2
# (right-aligned Arial, same as cells O2/Q2, etc.) instead of the old,
# slightly-different unused style it was carrying.
$ws.Range("O2").Copy()
$ws.Range("P2").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the placeholder that was missed previously.
$ws.Range("P2").Value = "{{renewal}}"
